# ---------------------------------------------------------------------------
# "thesis updated, assistant edit and delete functionality added"
#
# Translates the header row (and the "XB" size shorthand in column B) from
# English/abbreviated to Spanish, adds a thin border around the whole table,
# centers + wraps the table contents, widens several columns so the longer
# Spanish headers fit, fixes the row-1 height now that headers are shorter,
# and resets the sheet's scroll position back to A1.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Translate header row (row 1) to Spanish
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Modelo"
$ws.Range("B1").Value = "Cantidad de Parámetros"
$ws.Range("C1").Value = "Puntuación en tareas de Código"
$ws.Range("D1").Value = "Razonamiento común"
$ws.Range("E1").Value = "Conocimiento del Mundo"
$ws.Range("F1").Value = "Comprensión lectora"
$ws.Range("G1").Value = "Matemáticas"
# H1 (MMLU), I1 (BBH), J1 (AGI Eval) are unchanged acronyms.

# ---------------------------------------------------------------------------
# 2. Expand the "NB" (N billion) shorthand into "N Mil millones" for every
#    model-size cell in column B.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "7 Mil millones`n30 Mil millones"
$ws.Range("B3").Value = "7 Mil millones`n40 Mil millones"
$ws.Range("B4").Value = "7 Mil millones`n13 Mil millones`n33 Mil millones`n65 Mil millones"
$ws.Range("B5").Value = "7 Mil millones`n13 Mil millones`n34 Mil millones`n70 Mil millones"

# ---------------------------------------------------------------------------
# 3. Formatting: thin border around the whole table, centered + wrapped
#    text for the data area, and a border-only (no special alignment) style
#    for the model-name column (A).
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:J5")
$tableRange.Borders.Color = 0
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# Column A gets the border from $tableRange above; its alignment is left
# untouched (default/general - no centering, no wrap).

# Main data block (B1:J5): centered both ways + wrap text.
$dataRange = $ws.Range("B1:J5")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

# G1, I1, J1 (Matemáticas / BBH / AGI Eval): centered both ways, no wrap.
$noWrapHeaders = $ws.Range("G1,I1,J1")
$noWrapHeaders.HorizontalAlignment = -4108
$noWrapHeaders.VerticalAlignment = -4108
$noWrapHeaders.WrapText = $false

# C2 (MPT code score): centered horizontally + wrap, default (bottom)
# vertical alignment.
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").WrapText = $true
$ws.Range("C2").VerticalAlignment = -4107

# ---------------------------------------------------------------------------
# 4. Column widths - widen columns to fit the new, longer Spanish text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.7
$ws.Columns.Item(3).ColumnWidth = 15.85
$ws.Columns.Item(5).ColumnWidth = 12.7
$ws.Columns.Item(6).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 17.85

# ---------------------------------------------------------------------------
# 5. Row 1 is shorter now that the headers wrap onto fewer lines.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 6. Reset the view so the sheet opens scrolled to A1 (was topLeftCell=B1).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

Write-Host "Edit complete"
